$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Paragraph 1 (title): split "O que é Git e GitHub" into 3 runs, wrapping
# "Git" with proofErr spellStart/spellEnd. The trailing "?" run is untouched.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$body1 = '<w:p w14:paraId="293C64DE" w14:textId="77777777" w:rsidR="00721619" w:rsidRPr="00721619" w:rsidRDefault="00721619" w:rsidP="00721619">' +
    '<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="36"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:eastAsia="pt-BR"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="36"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:eastAsia="pt-BR"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve">O que é </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="36"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:eastAsia="pt-BR"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Git</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="36"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/><w:lang w:eastAsia="pt-BR"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve"> e GitHub</w:t></w:r>' +
    '<w:r w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="36"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:eastAsia="pt-BR"/><w14:ligatures w14:val="none"/></w:rPr><w:t>?</w:t></w:r>' +
    '</w:p>'
$p1.Range.InsertXML((New-PkgXml $body1)) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 2 ("Git : Chamado de repositório local, ..."): wrap "Git" with
# proofErr spellStart/gramStart...spellEnd, and close gramEnd after " :".
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$body2 = '<w:p w14:paraId="0E5A828B" w14:textId="4C757DC4" w:rsidR="001C0D56" w:rsidRDefault="0005095D">' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Git</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r w:rsidR="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> :</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidRPr="0005095D"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Chamado de repositório local</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>,</w:t></w:r>' +
    '<w:r w:rsidR="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidR="00721619" w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Usado para manter </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>vá</w:t></w:r>' +
    '<w:r w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ria</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s</w:t></w:r>' +
    '<w:r w:rsidR="00721619" w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> versões do software no computador</w:t></w:r>' +
    '<w:r w:rsidR="00E4193B"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r>' +
    '<w:r w:rsidR="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
$p2.Range.InsertXML((New-PkgXml $body2)) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 3 ("GitHub: Chamado de repositório remoto, ..."):
#  - pPr rPr sz/szCs 28 -> 24
#  - split "Usado para fazer..." into "Usado" (gramStart/gramEnd) + rest
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$body3 = '<w:p w14:paraId="3F32EB75" w14:textId="5266B227" w:rsidR="00721619" w:rsidRDefault="0005095D">' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="0005095D"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>GitHub</w:t></w:r>' +
    '<w:r w:rsidR="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Chamado de repositório remoto,</w:t></w:r>' +
    '<w:r w:rsidR="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r w:rsidR="00721619" w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Usado</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r w:rsidR="00721619" w:rsidRPr="00721619"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> para fazer o backup das versões na nuvem de forma gratuita</w:t></w:r>' +
    '<w:r w:rsidR="00E4193B"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r>' +
    '</w:p>'
$p3.Range.InsertXML((New-PkgXml $body3)) | Out-Null

# ---------------------------------------------------------------------------
# Insert two new paragraphs after paragraph 3 ("...gratuita.") and before the
# trailing empty paragraph: a centered "Para Gerar QR Code de um site"
# heading (with proofErr around "Code") and a link-style URL paragraph.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$insertPoint = $p4.Range
$insertPoint.Collapse(1)
$newParas = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Para Gerar QR </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Code</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> de um site</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>https://www.qr-code-generator.com/</w:t></w:r>' +
    '</w:p>'
$insertPoint.InsertXML((New-PkgXml $newParas)) | Out-Null
